$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 1076, pushing the existing data
# (previously rows 1076:1174) down to rows 1078:1176.
$ws.Rows.Item(1076).Insert()
$ws.Rows.Item(1076).Insert()

# New row 1076 (Calidad "Primera")
$ws.Range("A1076").Value = 3
$ws.Range("B1076").Value = "Femacal de La Calera"
$ws.Range("C1076").Value = "Coquimbo"
$ws.Range("D1076").Value = 45132
$ws.Range("E1076").Value = 5
$ws.Range("F1076").Value = 100114014
$ws.Range("G1076").Value = "Betarraga"
$ws.Range("H1076").Value = "Sin especificar"
$ws.Range("I1076").Value = "Primera"
$ws.Range("J1076").Value = 3750
$ws.Range("K1076").Value = 550
$ws.Range("L1076").Value = 600
$ws.Range("M1076").Value = 575
$ws.Range("N1076").Value = "$/paquete 4 unidades"
$ws.Range("O1076").Value = "Provincia de Quillota"
$ws.Range("P1076").Value = 144
$ws.Range("Q1076").Value = 4
$ws.Range("R1076").Value = "Hortaliza"

# New row 1077 (Calidad "Segunda")
$ws.Range("A1077").Value = 3
$ws.Range("B1077").Value = "Femacal de La Calera"
$ws.Range("C1077").Value = "Coquimbo"
$ws.Range("D1077").Value = 45132
$ws.Range("E1077").Value = 5
$ws.Range("F1077").Value = 100114014
$ws.Range("G1077").Value = "Betarraga"
$ws.Range("H1077").Value = "Sin especificar"
$ws.Range("I1077").Value = "Segunda"
$ws.Range("J1077").Value = 1600
$ws.Range("K1077").Value = 450
$ws.Range("L1077").Value = 450
$ws.Range("M1077").Value = 450
$ws.Range("N1077").Value = "$/paquete 4 unidades"
$ws.Range("O1077").Value = "Provincia de Quillota"
$ws.Range("P1077").Value = 112
$ws.Range("Q1077").Value = 4
$ws.Range("R1077").Value = "Hortaliza"

# Match the date format used by the rest of column D.
$ws.Range("D1076:D1077").NumberFormat = $ws.Range("D1078").NumberFormat
